# Adição do módulo de precificação
# Updates unit column ("Un" -> "Kg"/"g") for several existing products,
# fixes quantity/price for "file de peito de frango", and appends new
# product rows (16-20) for the pricing module.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows: unit changes from "Un" to "Kg" ---
$ws.Range("C4").Value  = "Kg"
$ws.Range("C5").Value  = "Kg"
$ws.Range("C7").Value  = "Kg"
$ws.Range("C8").Value  = "Kg"
$ws.Range("C11").Value = "Kg"
$ws.Range("C13").Value = "Kg"

# --- Update row 9 (file de peito de frango): quantity and price ---
$ws.Range("B9").Value = 1000
$ws.Range("D9").Value = 25

# --- Update row 15 (arroz parboilizado): unit and price ---
$ws.Range("C15").Value = "g"
$ws.Range("D15").Value = 3.5

# --- Append new rows for the pricing module ---
$ws.Range("A16").Value = "creme de leite"
$ws.Range("B16").Value = 5000
$ws.Range("C16").Value = "g"
$ws.Range("D16").Value = 30

$ws.Range("A17").Value = "batata palito mccain"
$ws.Range("B17").Value = 800
$ws.Range("C17").Value = "Kg"
$ws.Range("D17").Value = 34

$ws.Range("A18").Value = "queijo coalho"
$ws.Range("B18").Value = 0
$ws.Range("C18").Value = "g"
$ws.Range("D18").Value = 35

$ws.Range("A19").Value = "batata palha"
$ws.Range("B19").Value = 150
$ws.Range("C19").Value = "g"
$ws.Range("D19").Value = 30

$ws.Range("A20").Value = "cogumelo champignon"
$ws.Range("B20").Value = 300
$ws.Range("C20").Value = "g"
$ws.Range("D20").Value = 110
